$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update scenario labels (renumbered 1a/1b/2/3/4 -> 1/2/3/4/5, with minor wording tweaks)
$ws.Range("A10").Value = "1 – sea level rise of 2mm/year"
$ws.Range("A11").Value = "2 – As 1, with ntc of 0.15m"
$ws.Range("A12").Value = "3 – historic changes + slr"
$ws.Range("A13").Value = "4 – As 3, with dredge in 2000*"
$ws.Range("A14").Value = "5 – As 4, with reclamation in 2020*"

$ws.Range("A19").Value = "1 – sea level rise of 2mm/year"
$ws.Range("A20").Value = "2 – As 1, with ntc of 0.15m"
$ws.Range("A21").Value = "3 – historic changes + slr"
$ws.Range("A22").Value = "4 – As 3, with dredge in 2000*"
$ws.Range("A23").Value = "5 – As 4, with reclamation in 2020*"

# Update sheet view / selection: clear the frozen/topLeftCell view and select B42
$ws.Range("B42").Select()
